# [CSTPER-451] Import via XLS should offer stop in workspace, send to workflow and skip workflow
#
# Adds two new metadata columns (dc.type and dc.date.issued) to the "Main"
# sheet of the bulk-import test fixture, populates them for the three
# existing sample rows, and makes the "Main" sheet the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# New header cells for the two additional columns.
$ws.Range("D1").Value = "dc.type"
$ws.Range("E1").Value = "dc.date.issued"

# Row 2 - "My publication"
$ws.Range("D2").Value = "Article"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd"
$ws.Range("E2").Value = "2020-01-01"

# Row 3 - "Another publication"
$ws.Range("D3").Value = "Book"
$ws.Range("E3").NumberFormat = "yyyy-mm-dd"
$ws.Range("E3").Value = "2020-01-02"

# Row 4 - "Last publication"
$ws.Range("D4").Value = "Journal"
$ws.Range("E4").NumberFormat = "yyyy-mm-dd"
$ws.Range("E4").Value = "2020-01-03"

# An extra (empty) formatted cell below the data, matching the date format.
$ws.Range("E5").NumberFormat = "yyyy-mm-dd"

# Make "Main" the active sheet/tab, with the same selection left behind by
# the original edit.
$ws.Activate() | Out-Null
$ws.Range("E10").Select() | Out-Null
